$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "fd"
$ws.Range("C1").Value = "mnlk"
$ws.Range("C1").Select()
